$wb = $excel.ActiveWorkbook

# Insert the new "Drivers" sheet right after "Company_Profile" and before
# "ForgotPassword" (this also makes it the active/selected sheet, matching
# the target workbook view).
$afterSheet = $wb.Worksheets.Item("Company_Profile")
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Drivers"

# Column widths to match the rest of the workbook's "form" style sheets.
$ws.Columns.Item(1).ColumnWidth = 18
$ws.Columns.Item(2).ColumnWidth = 17.6640625

# Fill in the driver record. Values are entered in this order so that the
# generated shared-string table lines up with the source workbook.
$ws.Range("A1").Value = "First Name"
$ws.Range("A2").Value = "Last Name"
$ws.Range("B1").Value = "Abhijeet"
$ws.Range("B2").Value = "Nagarkar"
$ws.Range("A3").Value = "Cell"

$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:truckxdriver@gmail.com", "", "", "truckxdriver@gmail.com")
$ws.Range("B4").Style = "Hyperlink"

$ws.Range("A5").Value = "License"
$ws.Range("B3").Value = "800-793-9513"
$ws.Range("B5").Value = "b111111111111"
$ws.Range("A4").Value = "Email"

# The Device sheet had its selection changed to a full-row/column select
# (as if the user pressed Ctrl+A) rather than a single cell.
$device = $wb.Worksheets.Item("Device")
$device.Range("A1:XFD1048576").Select()
